# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Price column (D) cells being updated to Text format so that
# values such as "69.506.78" (thousand-separated, non-numeric) are preserved as
# literal text instead of being auto-parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.506.78'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '3.541.64'
$ws.Range("E3").Value = '  -1.75%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '195.81'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '583.02'
$ws.Range("E6").Value = '  -3.57%  '
$ws.Range("E7").Value = '  -2.42%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.205'
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("D11").Value = '51.79'
$ws.Range("E11").Value = '  -3.77%  '
$ws.Range("D12").Value = '0.0000287'
$ws.Range("E12").Value = '  -5.43%  '
$ws.Range("E13").Value = '  -3.73%  '
$ws.Range("D14").Value = '4.103.76'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").Value = '667.70'
$ws.Range("E15").Value = '  +12.04%  '
$ws.Range("D16").Value = '69.574.65'
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").Value = '12.55'
$ws.Range("E17").Value = '  -4.11%  '
$ws.Range("D18").Value = '3.538.88'
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("D20").Value = '18.45'
$ws.Range("E20").Value = '  -3.43%  '
$ws.Range("E21").Value = '  -3.22%  '
$ws.Range("D22").Value = '18.18'
$ws.Range("E22").Value = '  +2.34%  '
$ws.Range("D23").Value = '5.30'
$ws.Range("E23").Value = '  +2.55%  '
$ws.Range("D24").Value = '104.63'
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("E25").Value = '  -5.13%  '
$ws.Range("E26").Value = '  -4.04%  '
$ws.Range("D27").Value = '10.18'
$ws.Range("E27").Value = '  -5.24%  '
$ws.Range("D28").Value = '9.57'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '33.21'
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("D30").Value = '4.37'
$ws.Range("E30").Value = '  -7.72%  '
$ws.Range("D31").Value = '6.76'
$ws.Range("E31").Value = '  -5.50%  '
$ws.Range("D32").Value = '11.74'
$ws.Range("E32").Value = '  -4.31%  '
$ws.Range("E33").Value = '  -4.83%  '
$ws.Range("D34").Value = '61.81'
$ws.Range("E34").Value = '  -2.34%  '
$ws.Range("D35").Value = '3.784.31'
$ws.Range("E35").Value = '  -3.07%  '
$ws.Range("D36").Value = '0.0₃0813'
$ws.Range("E36").Value = '  -8.16%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '3.70'
$ws.Range("E38").Value = '  +4.74%  '
$ws.Range("D39").Value = '501.68'
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("E40").Value = '  -6.34%  '
$ws.Range("D41").Value = '0.371'
$ws.Range("E41").Value = '  -4.97%  '
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("D43").Value = '34.66'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("E47").Value = '  -2.76%  '
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '8.32'
$ws.Range("E49").Value = '  -3.58%  '
$ws.Range("D50").Value = '1.77'
$ws.Range("E50").Value = '  +19.01%  '
$ws.Range("D51").Value = '2.71'
$ws.Range("E51").Value = '  +63.49%  '

# Restore the original (default) cell style now that the text values are set,
# so the cells keep looking exactly as they did before (no explicit style index).
$ws.Range("D2:D51").Style = "Normal"

